$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = "FedEx® FedEx First Overnight®"
$ws.Range("B26").Value = "FedEx® FedEx Priority Overnight®"
$ws.Range("B27").Value = "FedEx® FedEx Standard Overnight®"
$ws.Range("B28").Value = "FedEx® FedEx 2Day® A.M."
$ws.Range("B29").Value = "FedEx® FedEx 2Day®"
$ws.Range("B30").Value = "FedEx® FedEx Express Saver®"

$ws.Range("B33").Select()
